$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear historical growth columns (D,E) for data rows 2 and 3 - removed in new schema
$ws.Range("D2:E3").ClearContents()

# Row 2 values (columns G..AQ)
$ws.Range("G2").Value = 0.08512820512820514
$ws.Range("H2").Value = 0.08512820512820514
$ws.Range("I2").Value = -0.02012820512820513
$ws.Range("J2").Value = -0.01771583597426294
$ws.Range("K2").Value = 0.235
$ws.Range("L2").Value = 0.03012820512820513
$ws.Range("M2").Value = 0.297
$ws.Range("N2").Value = 0.03421658986175115
$ws.Range("O2").Value = 1.263829787234042
$ws.Range("P2").Value = 0.297
$ws.Range("Q2").Value = 0.03421658986175115
$ws.Range("R2").Value = 1.263829787234042
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 1.97
$ws.Range("V2").Value = 0.2269585253456221
$ws.Range("W2").Value = 0.02790973871733967
$ws.Range("X2").Value = 0.06271493999845341
$ws.Range("Y2").Value = -0.03480520128111375
$ws.Range("Z2").Value = 1.695652173913044
$ws.Range("AA2").Value = -0.03003989578244586
$ws.Range("AB2").Value = 0.06271493999845341
$ws.Range("AC2").Value = -0.09275483578089927
$ws.Range("AD2").Value = 0
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 0
$ws.Range("AG2").Value = -1.97
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = -0.293591654247392
$ws.Range("AK2").Value = -0.309748427672956
$ws.Range("AL2").Value = 0.001
$ws.Range("AM2").Value = -0.415
$ws.Range("AN2").Value = -0
$ws.Range("AO2").Value = -157
$ws.Range("AP2").Value = 40.20408163265306
$ws.Range("AQ2").Value = 0.3783132530120482

# Row 3 values (columns G..AQ)
$ws.Range("G3").Value = 0.08512820512820514
$ws.Range("H3").Value = 0.08512820512820514
$ws.Range("I3").Value = -0.02012820512820513
$ws.Range("J3").Value = -0.01771583597426294
$ws.Range("K3").Value = 0.235
$ws.Range("L3").Value = 0.03012820512820513
$ws.Range("M3").Value = 0.297
$ws.Range("N3").Value = 0.03421658986175115
$ws.Range("O3").Value = 1.263829787234042
$ws.Range("P3").Value = 0.297
$ws.Range("Q3").Value = 0.03421658986175115
$ws.Range("R3").Value = 1.263829787234042
$ws.Range("S3").Value = 0
$ws.Range("T3").Value = 0
$ws.Range("U3").Value = 1.97
$ws.Range("V3").Value = 0.2269585253456221
$ws.Range("W3").Value = 0.02790973871733967
$ws.Range("X3").Value = 0.06271493999845341
$ws.Range("Y3").Value = -0.03480520128111375
$ws.Range("Z3").Value = 1.695652173913044
$ws.Range("AA3").Value = -0.03003989578244586
$ws.Range("AB3").Value = 0.06271493999845341
$ws.Range("AC3").Value = -0.09275483578089927
$ws.Range("AD3").Value = 0
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 0
$ws.Range("AG3").Value = -1.97
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = -0.293591654247392
$ws.Range("AK3").Value = -0.309748427672956
$ws.Range("AL3").Value = 0.001
$ws.Range("AM3").Value = -0.415
$ws.Range("AN3").Value = -0
$ws.Range("AO3").Value = -157
$ws.Range("AP3").Value = 40.20408163265306
$ws.Range("AQ3").Value = 0.3783132530120482
